$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.113.96'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '1.652.62'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5248'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.005'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2661'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06347'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.64'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07715'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.601'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("D13").Value = '1.679.99'
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").Value = '1.878.59'
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5619'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").Value = '0.0₅8191'
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.36'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.84%  '
$ws.Range("D18").Value = '26.119.81'
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.702'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.48%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.81%  '
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '191.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.985'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.006'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1204'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.263'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.92'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.508'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05631'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.272'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.494'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.371'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.579'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.796'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9507'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.411'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5751'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01595'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.989'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8386'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.76'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.04%  '
$ws.Range("D44").Value = '1.014.48'
$ws.Range("E44").Value = '  -5.76%  '
$ws.Range("D45").Value = '1.791.73'
$ws.Range("E45").Value = '  -0.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '58.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.65%  '
$ws.Range("E47").Value = '  +1.97%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.002'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.12%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05327'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4345'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.984'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.75%  '
